{"js": "// Apply the per-cell / header text replacements described by the diff:\n// the worksheet date line and each \"A\u00d7B=\" multiplication prompt are\n// swapped for a new value. Every old string in this document is unique,\n// so a simple search-and-replace (matchCase) per pair is sufficient and\n// avoids any row/column-index assumptions.\nconst replacements = [\n  [\"2025-03-26 Wednesday\", \"2025-03-27 Thursday\"],\n  [\"67\u00d738=\", \"23\u00d793=\"],\n  [\"98\u00d766=\", \"54\u00d764=\"],\n  [\"43\u00d763=\", \"59\u00d763=\"],\n  [\"95\u00d748=\", \"57\u00d768=\"],\n  [\"12\u00d774=\", \"92\u00d715=\"],\n  [\"15\u00d752=\", \"69\u00d753=\"],\n  [\"51\u00d732=\", \"46\u00d791=\"],\n  [\"15\u00d751=\", \"32\u00d735=\"],\n  [\"58\u00d752=\", \"45\u00d733=\"],\n  [\"86\u00d722=\", \"78\u00d767=\"],\n  [\"14\u00d725=\", \"69\u00d754=\"],\n  [\"88\u00d747=\", \"28\u00d722=\"],\n  [\"23\u00d748=\", \"77\u00d740=\"],\n  [\"44\u00d750=\", \"28\u00d741=\"],\n  [\"57\u00d779=\", \"31\u00d778=\"],\n  [\"59\u00d768=\", \"57\u00d785=\"],\n  [\"61\u00d747=\", \"19\u00d768=\"],\n  [\"92\u00d719=\", \"66\u00d759=\"],\n  [\"33\u00d722=\", \"79\u00d759=\"],\n  [\"60\u00d723=\", \"65\u00d739=\"],\n  [\"95\u00d795=\", \"34\u00d767=\"],\n  [\"98\u00d726=\", \"58\u00d717=\"],\n  [\"96\u00d774=\", \"92\u00d752=\"],\n  [\"73\u00d759=\", \"33\u00d783=\"],\n  [\"58\u00d751=\", \"81\u00d770=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the per-cell / header text replacements described by the diff:\n# the worksheet date line and each \"A\u00d7B=\" multiplication prompt are\n# swapped for a new value. Every \"old\" string in this document is\n# unique, so a plain Find/Replace-all pass per pair is sufficient and\n# avoids any row/column-index assumptions.\n\n$replacements = @(\n    @{ Find = \"2025-03-26 Wednesday\"; Replace = \"2025-03-27 Thursday\" },\n    @{ Find = \"67\u00d738=\"; Replace = \"23\u00d793=\" },\n    @{ Find = \"98\u00d766=\"; Replace = \"54\u00d764=\" },\n    @{ Find = \"43\u00d763=\"; Replace = \"59\u00d763=\" },\n    @{ Find = \"95\u00d748=\"; Replace = \"57\u00d768=\" },\n    @{ Find = \"12\u00d774=\"; Replace = \"92\u00d715=\" },\n    @{ Find = \"15\u00d752=\"; Replace = \"69\u00d753=\" },\n    @{ Find = \"51\u00d732=\"; Replace = \"46\u00d791=\" },\n    @{ Find = \"15\u00d751=\"; Replace = \"32\u00d735=\" },\n    @{ Find = \"58\u00d752=\"; Replace = \"45\u00d733=\" },\n    @{ Find = \"86\u00d722=\"; Replace = \"78\u00d767=\" },\n    @{ Find = \"14\u00d725=\"; Replace = \"69\u00d754=\" },\n    @{ Find = \"88\u00d747=\"; Replace = \"28\u00d722=\" },\n    @{ Find = \"23\u00d748=\"; Replace = \"77\u00d740=\" },\n    @{ Find = \"44\u00d750=\"; Replace = \"28\u00d741=\" },\n    @{ Find = \"57\u00d779=\"; Replace = \"31\u00d778=\" },\n    @{ Find = \"59\u00d768=\"; Replace = \"57\u00d785=\" },\n    @{ Find = \"61\u00d747=\"; Replace = \"19\u00d768=\" },\n    @{ Find = \"92\u00d719=\"; Replace = \"66\u00d759=\" },\n    @{ Find = \"33\u00d722=\"; Replace = \"79\u00d759=\" },\n    @{ Find = \"60\u00d723=\"; Replace = \"65\u00d739=\" },\n    @{ Find = \"95\u00d795=\"; Replace = \"34\u00d767=\" },\n    @{ Find = \"98\u00d726=\"; Replace = \"58\u00d717=\" },\n    @{ Find = \"96\u00d774=\"; Replace = \"92\u00d752=\" },\n    @{ Find = \"73\u00d759=\"; Replace = \"33\u00d783=\" },\n    @{ Find = \"58\u00d751=\"; Replace = \"81\u00d770=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # 0 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 0, $false, $pair.Replace, 2)\n}\n"}
